$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage each new value as forced-text before a
# values-only paste into the target cell. This avoids Excel silently
# re-interpreting numeric-looking strings (e.g. "519.87", "1.00") as
# numbers, while still leaving the destination cell style untouched
# (PasteSpecial xlPasteValues copies only the value, not formatting).
$tmp = $ws.Range("Z1")
$tmp.NumberFormat = "@"

$tmp.Value = "57.012.87"
$tmp.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$tmp.Value = "  -1.46%  "
$tmp.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$tmp.Value = "3.083.09"
$tmp.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$tmp.Value = "  -0.22%  "
$tmp.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$tmp.Value = "  +0.06%  "
$tmp.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$tmp.Value = "519.87"
$tmp.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$tmp.Value = "  -1.20%  "
$tmp.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$tmp.Value = "135.83"
$tmp.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$tmp.Value = "  -3.58%  "
$tmp.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$tmp.Value = "3.081.87"
$tmp.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$tmp.Value = "  -0.21%  "
$tmp.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$tmp.Value = "0.452"
$tmp.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$tmp.Value = "  +2.28%  "
$tmp.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$tmp.Value = "  +2.59%  "
$tmp.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$tmp.Value = "  -1.48%  "
$tmp.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$tmp.Value = "0.398"
$tmp.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$tmp.Value = "  +1.60%  "
$tmp.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$tmp.Value = "3.615.05"
$tmp.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$tmp.Value = "  -0.09%  "
$tmp.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$tmp.Value = "  +1.72%  "
$tmp.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$tmp.Value = "25.28"
$tmp.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$tmp.Value = "  -0.62%  "
$tmp.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$tmp.Value = "  -2.06%  "
$tmp.Copy()
$ws.Range("E16").PasteSpecial(-4163)
$tmp.Value = "57.089.42"
$tmp.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$tmp.Value = "  -1.31%  "
$tmp.Copy()
$ws.Range("E17").PasteSpecial(-4163)
$tmp.Value = "3.082.77"
$tmp.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$tmp.Value = "  +0.17%  "
$tmp.Copy()
$ws.Range("E18").PasteSpecial(-4163)
$tmp.Value = "  -3.17%  "
$tmp.Copy()
$ws.Range("E19").PasteSpecial(-4163)
$tmp.Value = "  -1.60%  "
$tmp.Copy()
$ws.Range("E20").PasteSpecial(-4163)
$tmp.Value = "  -1.29%  "
$tmp.Copy()
$ws.Range("E21").PasteSpecial(-4163)
$tmp.Value = "346.56"
$tmp.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$tmp.Value = "  +1.25%  "
$tmp.Copy()
$ws.Range("E22").PasteSpecial(-4163)
$tmp.Value = "Dai"
$tmp.Copy()
$ws.Range("B23").PasteSpecial(-4163)
$tmp.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$tmp.Copy()
$ws.Range("C23").PasteSpecial(-4163)
$tmp.Value = "1.00"
$tmp.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$tmp.Value = "  +0.02%  "
$tmp.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$tmp.Value = "LEO"
$tmp.Copy()
$ws.Range("B24").PasteSpecial(-4163)
$tmp.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$tmp.Copy()
$ws.Range("C24").PasteSpecial(-4163)
$tmp.Value = "5.77"
$tmp.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$tmp.Value = "  +1.24%  "
$tmp.Copy()
$ws.Range("E24").PasteSpecial(-4163)
$tmp.Value = "68.16"
$tmp.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$tmp.Value = "  +1.29%  "
$tmp.Copy()
$ws.Range("E25").PasteSpecial(-4163)
$tmp.Value = "  -2.44%  "
$tmp.Copy()
$ws.Range("E26").PasteSpecial(-4163)
$tmp.Value = "  -1.22%  "
$tmp.Copy()
$ws.Range("E27").PasteSpecial(-4163)
$tmp.Value = "1.00"
$tmp.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$tmp.Value = "  +0.17%  "
$tmp.Copy()
$ws.Range("E28").PasteSpecial(-4163)
$tmp.Value = "0.0₃0861"
$tmp.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$tmp.Value = "  -5.75%  "
$tmp.Copy()
$ws.Range("E29").PasteSpecial(-4163)
$tmp.Value = "7.25"
$tmp.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$tmp.Value = "  +0.17%  "
$tmp.Copy()
$ws.Range("E31").PasteSpecial(-4163)
$tmp.Value = "  -0.27%  "
$tmp.Copy()
$ws.Range("E32").PasteSpecial(-4163)
$tmp.Value = "  -8.64%  "
$tmp.Copy()
$ws.Range("E33").PasteSpecial(-4163)
$tmp.Value = "20.80"
$tmp.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$tmp.Value = "  -0.61%  "
$tmp.Copy()
$ws.Range("E34").PasteSpecial(-4163)
$tmp.Value = "  +6.47%  "
$tmp.Copy()
$ws.Range("E35").PasteSpecial(-4163)
$tmp.Value = "159.32"
$tmp.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$tmp.Value = "  +0.52%  "
$tmp.Copy()
$ws.Range("E36").PasteSpecial(-4163)
$tmp.Value = "1.13"
$tmp.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$tmp.Value = "  -4.07%  "
$tmp.Copy()
$ws.Range("E37").PasteSpecial(-4163)
$tmp.Value = "5.99"
$tmp.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$tmp.Value = "  -2.39%  "
$tmp.Copy()
$ws.Range("E38").PasteSpecial(-4163)
$tmp.Value = "25.71"
$tmp.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$tmp.Value = "  -0.68%  "
$tmp.Copy()
$ws.Range("E39").PasteSpecial(-4163)
$tmp.Value = "  -0.36%  "
$tmp.Copy()
$ws.Range("E40").PasteSpecial(-4163)
$tmp.Value = "  -2.14%  "
$tmp.Copy()
$ws.Range("E41").PasteSpecial(-4163)
$tmp.Value = "  +1.45%  "
$tmp.Copy()
$ws.Range("E42").PasteSpecial(-4163)
$tmp.Value = "  +0.35%  "
$tmp.Copy()
$ws.Range("E43").PasteSpecial(-4163)
$tmp.Value = "  +0.99%  "
$tmp.Copy()
$ws.Range("E44").PasteSpecial(-4163)
$tmp.Value = "2.392.92"
$tmp.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$tmp.Value = "  +5.06%  "
$tmp.Copy()
$ws.Range("E45").PasteSpecial(-4163)
$tmp.Value = "36.62"
$tmp.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$tmp.Value = "  -0.79%  "
$tmp.Copy()
$ws.Range("E46").PasteSpecial(-4163)
$tmp.Value = "  +0.08%  "
$tmp.Copy()
$ws.Range("E47").PasteSpecial(-4163)
$tmp.Value = "3.123.47"
$tmp.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$tmp.Value = "  -0.14%  "
$tmp.Copy()
$ws.Range("E48").PasteSpecial(-4163)
$tmp.Value = "0.0264"
$tmp.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$tmp.Value = "  +0.86%  "
$tmp.Copy()
$ws.Range("E49").PasteSpecial(-4163)
$tmp.Value = "0.958"
$tmp.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$tmp.Value = "  -3.46%  "
$tmp.Copy()
$ws.Range("E50").PasteSpecial(-4163)
$tmp.Value = "  -2.38%  "
$tmp.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$tmp.Clear()
$excel.CutCopyMode = 0
